# Actualizacion Datos Personales 4 nov
# Inserts a new student record into the "Rescatables" sheet (row 4),
# shifting the existing rows 4-7 down to rows 5-8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Insert a new row before row 4 so the following rows shift down.
$ws.Rows.Item(4).Insert()

# Fill the new row 4 with the new student's data.
$ws.Cells.Item(4, 1).Value = 19330051920235
$ws.Cells.Item(4, 2).Value = "HERAS"
$ws.Cells.Item(4, 3).Value = "LOPEZ"
$ws.Cells.Item(4, 4).Value = "CESAR ENRIQUE"
$ws.Cells.Item(4, 5).Value = "CIENCIA, TECNOLOGÍA, SOCIEDAD Y VALORES"
$ws.Cells.Item(4, 6).Value = "5APM"
$ws.Cells.Item(4, 7).Value = 6
